# Ajustes importantes para o correto funcionamento do BD
#
# 1) The empty paragraph right after the "2)Instruções para execução" heading
#    gets a new explanatory paragraph of text (with the "posições." grammar
#    proofing marks preserved) and the "_GoBack" bookmark moves to the end
#    of that new text.
# 2) The old paragraph at the end of the document that used to hold the
#    "_GoBack" bookmark becomes an empty paragraph again.

$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark wherever Word last left it ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the heading paragraph "2)Instruções para execução" ---
$found = $d.Content
$found.Find.Execute("2)Instruções para execução", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($found.Start -ge $p.Range.Start -and $found.Start -lt $p.Range.End) {
        $headingIndex = $i
        break
    }
}

# --- Step 3: fill in the first empty paragraph right after the heading ---
$target = $d.Paragraphs.Item($headingIndex + 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Toda vez que a aplicação estiver sendo iniciada, todos os dados de posição que são relacionados aos veículos que foram fornecidos na planilha “</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>posições.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>csv” serão excluídos e inseridos novamente na base de dados com informações idênticas as contidas na planilha em questão.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)
